# saved_addresses_2.xlsx -- "upload funkcniho vraceni zmen"
#
# Restores/rewrites the address rows across ip_address_list,
# ip_adress_fav_list and projects_bin2, flips projects_bin2 back to
# visible and makes it the active sheet, and removes the unused
# "Hyperlink" cell style that nothing in the workbook references anymore.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: ip_address_list
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("ip_address_list")

$ws1.Range("A1").Value = "aaaaaaaa"
$ws1.Range("B1").Value = "192.168.1.131"
$ws1.Range("D1").Value = "ccxxxggd"

$ws1.Range("A2").Value = "Domaci Wifiaffz"
$ws1.Range("B2").Value = "192.168.1.13¨ks"
$ws1.Range("D2").Value = "ddassssaa"
$ws1.Range("E2").Value = $false

$ws1.Range("A3").Value = "514nnnzzzzz"
$ws1.Range("B3").Value = "192.168.14.240a"
$ws1.Range("D3").ClearContents()
$ws1.Range("E3").Value = 1

# ---------------------------------------------------------------------
# Sheet 2: ip_adress_fav_list
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("ip_adress_fav_list")

$ws2.Range("A1").Value = "514_Teleflex"
$ws2.Range("B1").Value = "192.168.14.240a"

$ws2.Range("A2").Value = "514nnnzzzzz"
$ws2.Range("B2").Value = "192.168.14.240a"
$ws2.Range("E2").Value = 1

# ---------------------------------------------------------------------
# Sheet 5: projects_bin2 -- unhide it, rewrite its single data row and
# shift that row from row 1 down to row 2.
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("projects_bin2")
$ws5.Visible = -1

$ws5.Range("A1:E1").ClearContents()
$ws5.Range("A2").Value = "Domaci Wifiaffz"
$ws5.Range("B2").Value = "192.168.1.13¨ks"
$ws5.Range("C2").Value = "255.255.255.0"
$ws5.Range("D2").Value = "ddassssaa"
$ws5.Range("E2").Value = $false

# Make projects_bin2 the active sheet/tab (also updates workbookView's
# activeTab and clears tabSelected from whichever sheet had it before).
$ws5.Activate()

# ---------------------------------------------------------------------
# Styles: drop the now-unused "Hyperlink" cell style.
# ---------------------------------------------------------------------
$wb.Styles.Item("Hyperlink").Delete()
